# Update the division expressions in the practice-worksheet table cells.
# Each call performs an exact, whole-document "replace all" of a unique
# "a÷b=" expression with its new value (MatchWildcards = $true is used so
# the literal text is matched verbatim; wdReplaceAll = 2 replaces every
# occurrence, but since each "old" string occurs exactly once in the
# document this is equivalent to a single targeted substitution).
$d = $word.ActiveDocument

$d.Content.Find.Execute("14÷5=", $false, $false, $false, $false, $false, $true, 1, $false, "10÷3=", 2)
$d.Content.Find.Execute("93÷5=", $false, $false, $false, $false, $false, $true, 1, $false, "63÷8=", 2)
$d.Content.Find.Execute("42÷3=", $false, $false, $false, $false, $false, $true, 1, $false, "34÷2=", 2)
$d.Content.Find.Execute("96÷8=", $false, $false, $false, $false, $false, $true, 1, $false, "46÷3=", 2)
$d.Content.Find.Execute("12÷7=", $false, $false, $false, $false, $false, $true, 1, $false, "78÷9=", 2)
$d.Content.Find.Execute("61÷6=", $false, $false, $false, $false, $false, $true, 1, $false, "56÷5=", 2)
$d.Content.Find.Execute("56÷4=", $false, $false, $false, $false, $false, $true, 1, $false, "56÷7=", 2)
$d.Content.Find.Execute("60÷8=", $false, $false, $false, $false, $false, $true, 1, $false, "32÷6=", 2)
$d.Content.Find.Execute("20÷5=", $false, $false, $false, $false, $false, $true, 1, $false, "19÷4=", 2)
$d.Content.Find.Execute("42÷7=", $false, $false, $false, $false, $false, $true, 1, $false, "25÷6=", 2)
$d.Content.Find.Execute("57÷8=", $false, $false, $false, $false, $false, $true, 1, $false, "46÷7=", 2)
$d.Content.Find.Execute("70÷4=", $false, $false, $false, $false, $false, $true, 1, $false, "33÷6=", 2)
$d.Content.Find.Execute("56÷6=", $false, $false, $false, $false, $false, $true, 1, $false, "73÷4=", 2)
$d.Content.Find.Execute("96÷5=", $false, $false, $false, $false, $false, $true, 1, $false, "86÷3=", 2)
$d.Content.Find.Execute("67÷3=", $false, $false, $false, $false, $false, $true, 1, $false, "94÷2=", 2)
$d.Content.Find.Execute("58÷7=", $false, $false, $false, $false, $false, $true, 1, $false, "67÷3=", 2)
$d.Content.Find.Execute("50÷7=", $false, $false, $false, $false, $false, $true, 1, $false, "61÷8=", 2)
$d.Content.Find.Execute("29÷5=", $false, $false, $false, $false, $false, $true, 1, $false, "43÷2=", 2)
$d.Content.Find.Execute("47÷9=", $false, $false, $false, $false, $false, $true, 1, $false, "42÷6=", 2)
$d.Content.Find.Execute("46÷6=", $false, $false, $false, $false, $false, $true, 1, $false, "49÷8=", 2)
$d.Content.Find.Execute("44÷5=", $false, $false, $false, $false, $false, $true, 1, $false, "17÷6=", 2)
$d.Content.Find.Execute("60÷4=", $false, $false, $false, $false, $false, $true, 1, $false, "21÷6=", 2)
$d.Content.Find.Execute("34÷3=", $false, $false, $false, $false, $false, $true, 1, $false, "41÷6=", 2)
$d.Content.Find.Execute("49÷4=", $false, $false, $false, $false, $false, $true, 1, $false, "36÷6=", 2)
$d.Content.Find.Execute("75÷6=", $false, $false, $false, $false, $false, $true, 1, $false, "80÷9=", 2)
